$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

$sh = $grp.GroupItems.Item("tx12")
$sh.Left = 453.3802
$sh.Top = 353.55141

$sh = $grp.GroupItems.Item("tx13")
$sh.Left = 466.1839
$sh.Top = 381.5106

$sh = $grp.GroupItems.Item("tx14")
$sh.Left = 318.7779
$sh.Top = 390.0324

$sh = $grp.GroupItems.Item("tx15")
$sh.Left = 337.7418
$sh.Top = 417.7998

$sh = $grp.GroupItems.Item("tx16")
$sh.Left = 309.313
$sh.Top = 334.8286

$sh = $grp.GroupItems.Item("tx17")
$sh.Left = 329.2437
$sh.Top = 362.58773

$sh = $grp.GroupItems.Item("tx18")
$sh.Left = 262.9012
$sh.Top = 291.08962

$sh = $grp.GroupItems.Item("tx19")
$sh.Left = 277.6014
$sh.Top = 315.4561

$sh = $grp.GroupItems.Item("tx20")
$sh.Left = 341.0384
$sh.Top = 293.239

$sh = $grp.GroupItems.Item("tx21")
$sh.Left = 364.7578
$sh.Top = 317.6055

$sh = $grp.GroupItems.Item("tx22")
$sh.Left = 288.57961
$sh.Top = 246.9846

$sh = $grp.GroupItems.Item("tx23")
$sh.Left = 293.3144
$sh.Top = 271.35103

$sh = $grp.GroupItems.Item("tx24")
$sh.Left = 327.6419
$sh.Top = 199.4468

$sh = $grp.GroupItems.Item("tx25")
$sh.Left = 347.56441
$sh.Top = 227.4143

$sh = $grp.GroupItems.Item("tx26")
$sh.Left = 431.6339
$sh.Top = 193.4866

$sh = $grp.GroupItems.Item("tx27")
$sh.Left = 440.6448
$sh.Top = 217.8447
